# Adjust Investment Summary table column widths for better formatting
#
# Three placeholder tables (slides 2, 3, 4) each get their last grid
# column widened by 1 EMU so the overall table/frame extent becomes
# 8710933 EMU (instead of 8710932 EMU), and all of the placeholder
# sample text that was typed into the cells is cleared back out.
#
# 1 EMU = 1/12700 pt, so widen the final column by exactly 1 EMU worth
# of points: (old_emu + 1) / 12700.

$p = $ppt.ActivePresentation

function Clear-TableCells($tbl) {
    for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
        for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
            $tbl.Cell($r, $c).Shape.TextFrame.TextRange.Text = ""
        }
    }
}

# --- Slide 2: "Why This Solution?" (2 columns x 4 rows) ---
$s2 = $p.Slides.Item(2)
$tbl2 = $s2.Shapes.Item(3).Table
$tbl2.Columns.Item(2).Width = 4355467 / 12700.0
Clear-TableCells $tbl2

# --- Slide 3: "Business Value - Financial Impact" (2 columns x 6 rows) ---
$s3 = $p.Slides.Item(3)
$tbl3 = $s3.Shapes.Item(3).Table
$tbl3.Columns.Item(2).Width = 4355467 / 12700.0
Clear-TableCells $tbl3

# --- Slide 4: "Risk Mitigation" (3 columns x 4 rows) ---
$s4 = $p.Slides.Item(4)
$tbl4 = $s4.Shapes.Item(3).Table
$tbl4.Columns.Item(3).Width = 2903645 / 12700.0
Clear-TableCells $tbl4
